# tiretrainingdata.xlsx — "Kvals not needed in mu calculation, hence removed"
#
# Each of the three 12-row data blocks (A3:C14, A16:C27, A29:C40) was
# previously scaled by a per-block "K" constant. That constant equals the
# block's own last B-column value (B14 / B27 / B40), since after the edit
# each block's last B value becomes exactly 1. Removing the K factor means
# dividing every cell in the block (columns A:C) by that original B value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$blocks = @(
    @{ Start = 3;  End = 14 },
    @{ Start = 16; End = 27 },
    @{ Start = 29; End = 40 }
)

foreach ($block in $blocks) {
    $start = $block.Start
    $end   = $block.End

    # K value for this block = the block's own last B-column entry
    $kval = $ws.Cells.Item($end, 2).Value2

    $rng  = $ws.Range($ws.Cells.Item($start, 1), $ws.Cells.Item($end, 3))
    $vals = $rng.Value2

    for ($i = 1; $i -le $vals.GetLength(0); $i++) {
        for ($j = 1; $j -le $vals.GetLength(1); $j++) {
            $vals[$i, $j] = $vals[$i, $j] / $kval
        }
    }

    $rng.Value2 = $vals
}

# Selection moved from the whole-data range to the single cell D8
$ws.Range("D8").Select()
